$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'94.379.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.53%  "

$ws.Range("D3").Value = "'3.423.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.58%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'237.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.95%  "

$ws.Range("D6").Value = "'642.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.32%  "

$ws.Range("D7").Value = "'1.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").Value = "'0.405"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.84%  "

$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").Value = "'0.970"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.51%  "

$ws.Range("D11").Value = "'3.420.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.56%  "

$ws.Range("E12").Value = "  -4.61%  "

$ws.Range("D13").Value = "'41.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.25%  "

$ws.Range("E14").Value = "  +2.14%  "

$ws.Range("D15").Value = "'94.148.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.38%  "

$ws.Range("D16").Value = "'4.066.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.68%  "

$ws.Range("E17").Value = "  -0.96%  "

$ws.Range("D18").Value = "'8.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.22%  "

$ws.Range("D19").Value = "'3.427.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.81%  "

$ws.Range("D20").Value = "'17.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.88%  "

$ws.Range("E21").Value = "  +6.15%  "

$ws.Range("D22").Value = "'0.498"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.00%  "

$ws.Range("D23").Value = "'497.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.28%  "

$ws.Range("E24").Value = "  -5.25%  "

$ws.Range("E25").Value = "  -3.23%  "

$ws.Range("D26").Value = "'6.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.26%  "

$ws.Range("D27").Value = "'93.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.57%  "

$ws.Range("D28").Value = "'11.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.68%  "

$ws.Range("D29").Value = "'3.607.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.68%  "

$ws.Range("D30").Value = "'11.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.94%  "

$ws.Range("E31").Value = "  -0.21%  "

$ws.Range("E32").Value = "  +8.51%  "

$ws.Range("E33").Value = "  -1.71%  "

$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "

$ws.Range("E35").Value = "  -4.01%  "

$ws.Range("D36").Value = "'29.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.53%  "

$ws.Range("D37").Value = "'0.552"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("D38").Value = "'548.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.57%  "

$ws.Range("D39").Value = "'7.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.15%  "

$ws.Range("E40").Value = "  -2.57%  "

$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("D43").Value = "'0.900"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.01%  "

$ws.Range("D44").Value = "'24.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.33%  "

$ws.Range("E45").Value = "  -0.25%  "

$ws.Range("B46").Value = "MantraDAO"
$ws.Range("C46").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D46").Value = "'3.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.51%  "

$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'5.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "

$ws.Range("D48").Value = "'3.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.11%  "

$ws.Range("D49").Value = "'0.0409"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.35%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.60%  "

$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "'54.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.24%  "
